$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 0.3491903333333333
$ws.Range("H2").Value = 1.047571
$ws.Range("I2").Value = 0.008130334326258625
$ws.Range("J2").Value = 0.008130334326258625
$ws.Range("M2").Value = 7.413580666666667
$ws.Range("N2").Value = 22.240742
$ws.Range("O2").Value = 0.05108888817597561
$ws.Range("P2").Value = 0.05108888817597561
$ws.Range("Q2").Value = 2.588750704186889
$ws.Range("R2").Value = 23.298756337682
$ws.Range("S2").Value = 0.000415369741227523
$ws.Range("T2").Value = 0.0004153697412275229
$ws.Range("G3").Value = 0.3491903333333333
$ws.Range("H3").Value = 1.047571
$ws.Range("I3").Value = 0.008130334326258625
$ws.Range("J3").Value = 0.008130334326258625
$ws.Range("O3").Value = 0.0112127179963522
$ws.Range("P3").Value = 0.0112127179963522
$ws.Range("Q3").Value = 0.5681652634311112
$ws.Range("R3").Value = 5.113487370880001
$ws.Range("S3").Value = 0.00009116314601640014
$ws.Range("T3").Value = 0.00009116314601640014
$ws.Range("G4").Value = 0.3491903333333333
$ws.Range("H4").Value = 1.047571
$ws.Range("I4").Value = 0.008130334326258625
$ws.Range("J4").Value = 0.008130334326258625
$ws.Range("O4").Value = 0.9376983938276722
$ws.Range("P4").Value = 0.9376983938276722
$ws.Range("Q4").Value = 47.51458612633911
$ws.Range("R4").Value = 427.6312751370521
$ws.Range("S4").Value = 0.007623801439014703
$ws.Range("T4").Value = 0.007623801439014703
$ws.Range("I5").Value = 0.801301577139928
$ws.Range("J5").Value = 0.8013015771399279
$ws.Range("M5").Value = 7.413580666666667
$ws.Range("N5").Value = 22.240742
$ws.Range("O5").Value = 0.05108888817597561
$ws.Range("P5").Value = 0.05108888817597561
$ws.Range("Q5").Value = 255.1395722298207
$ws.Range("R5").Value = 2296.256150068386
$ws.Range("S5").Value = 0.04093760666973468
$ws.Range("T5").Value = 0.04093760666973467
$ws.Range("I6").Value = 0.801301577139928
$ws.Range("J6").Value = 0.8013015771399279
$ws.Range("O6").Value = 0.0112127179963522
$ws.Range("P6").Value = 0.0112127179963522
$ws.Range("S6").Value = 0.008984768614502273
$ws.Range("T6").Value = 0.008984768614502271
$ws.Range("I7").Value = 0.801301577139928
$ws.Range("J7").Value = 0.8013015771399279
$ws.Range("O7").Value = 0.9376983938276722
$ws.Range("P7").Value = 0.9376983938276722
$ws.Range("S7").Value = 0.7513792018556911
$ws.Range("T7").Value = 0.751379201855691
$ws.Range("I8").Value = 0.1905680885338134
$ws.Range("J8").Value = 0.1905680885338134
$ws.Range("M8").Value = 7.413580666666667
$ws.Range("N8").Value = 22.240742
$ws.Range("O8").Value = 0.05108888817597561
$ws.Range("P8").Value = 0.05108888817597561
$ws.Range("Q8").Value = 60.67810419482201
$ws.Range("R8").Value = 546.102937753398
$ws.Range("S8").Value = 0.009735911765013414
$ws.Range("T8").Value = 0.009735911765013412
$ws.Range("I9").Value = 0.1905680885338134
$ws.Range("J9").Value = 0.1905680885338134
$ws.Range("O9").Value = 0.0112127179963522
$ws.Range("P9").Value = 0.0112127179963522
$ws.Range("S9").Value = 0.002136786235833529
$ws.Range("T9").Value = 0.002136786235833529
$ws.Range("I10").Value = 0.1905680885338134
$ws.Range("J10").Value = 0.1905680885338134
$ws.Range("O10").Value = 0.9376983938276722
$ws.Range("P10").Value = 0.9376983938276722
$ws.Range("S10").Value = 0.1786953905329665
$ws.Range("T10").Value = 0.1786953905329665
